$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Bus m -> id -> r -> x -> b -> Rating -> Costo)
# This naturally inherits the formatting of the column to its left (column C),
# matching the styles used by the rest of the header/data rows.
$ws.Columns("D:D").Insert()

# New header cell for the inserted "id" column
$ws.Range("D1").Value = "id"

# New data cell value for row 2
$ws.Range("D2").Value = "NL"

# Update selection to match the authored state
$ws.Range("D3").Select()
